$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status fields for row 2: ID_STATUS -> 2, STATUS -> "Inactivo(a)"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = "Inactivo(a)"

# Reflect the last-edited/selected cell as in the authored workbook
$ws.Range("F4").Select()
